# Modify temperature data form
# - "survey" sheet rows 9-15: replace the four temperature-measurement fields
#   (average/lower-alarm/minimum/cumulative-below, upper-alarm/maximum/cumulative-above)
#   with four new 30-day alarm/threshold count fields, and blank out the rows
#   that are no longer needed.
# - Update the active sheet/selection state to reflect the sheet last worked on.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# type column - all four new fields are integer counts/day-counts
$ws.Range("C9").Value = "integer"
$ws.Range("C10").Value = "integer"
$ws.Range("C11").Value = "integer"
$ws.Range("C12").Value = "integer"

# name column
$ws.Range("E9").Value = "number_of_high_alarms_30"
$ws.Range("E10").Value = "number_of_low_alarms_30"
$ws.Range("E11").Value = "days_temp_above_8_30"
$ws.Range("E12").Value = "days_temp_below_2_30"

# display.prompt.text (English label) column
$ws.Range("F9").Value = "Number of High Alarms Over Last 30 Days"
$ws.Range("F10").Value = "Number of Low Alarms Over Last 30 Days"
$ws.Range("F11").Value = "Days with Temperatures Above 8°C, Last 30 Days"
$ws.Range("F12").Value = "Days with Temperature Below 2°C, Last 30 Days"

# display.prompt.text.es (Spanish label) column
$ws.Range("G9").Value = "Número de alarmas altas en los últimos 30 días"
$ws.Range("G10").Value = "Número de alarmas bajas en los últimos 30 días"
$ws.Range("G11").Value = "Días con temperaturas superiores a 8 ° C, últimos 30 días"
$ws.Range("G12").Value = "Días con temperatura inferior a 2 ° C, últimos 30 días"

# display.hint.text / display.hint.text.es (units) no longer apply - clear them
$ws.Range("I9").Value = ""
$ws.Range("J9").Value = ""
$ws.Range("I10").Value = ""
$ws.Range("I11").Value = ""
$ws.Range("J11").Value = ""
$ws.Range("I12").Value = ""
$ws.Range("J12").Value = ""

# Rows 13-15 no longer hold a field definition - clear their contents.
$ws.Range("C13").Value = ""
$ws.Range("E13").Value = ""
$ws.Range("F13").Value = ""
$ws.Range("G13").Value = ""
$ws.Range("I13").Value = ""

$ws.Range("C14").Value = ""
$ws.Range("E14").Value = ""
$ws.Range("F14").Value = ""
$ws.Range("G14").Value = ""
$ws.Range("I14").Value = ""
$ws.Range("J14").Value = ""

$ws.Range("C15").Value = ""
$ws.Range("E15").Value = ""
$ws.Range("F15").Value = ""
$ws.Range("G15").Value = ""
$ws.Range("I15").Value = ""
$ws.Range("J15").Value = ""

# Row 9 wraps the (now longer) Spanish label onto two lines.
$ws.Range("G9").WrapText = $true
$ws.Rows.Item(9).RowHeight = 31.2

# Reflect the final selection/active-sheet state: "settings" loses focus,
# "survey" (the sheet being edited) becomes the active tab/selection.
$settings = $wb.Worksheets.Item("settings")
$settings.Activate()
$settings.Range("F8").Select()

$ws.Activate()
$ws.Range("G15").Select()
